$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 3603
$wsExpo.Range("F5").Value = 8390
$wsExpo.Range("F8").Value = 2279
$wsExpo.Range("F11").Value = 680
$wsExpo.Range("F13").Value = 7616
$wsExpo.Range("F14").Value = 7828
$wsExpo.Range("F15").Value = 0
$wsExpo.Range("F16").Value = 4980
$wsExpo.Range("F18").Value = 539
$wsExpo.Range("G20").Value = 58.8
$wsExpo.Range("F23").Value = 5355
$wsExpo.Range("G27").Value = "不可售"
$wsExpo.Range("F28").Value = 950
$wsExpo.Range("F29").Value = 1465
$wsExpo.Range("F30").Value = 2071
$wsExpo.Range("F31").Value = 26
$wsExpo.Range("F32").Value = 198
$wsExpo.Range("F34").Value = 1099
$wsExpo.Range("F35").Value = 9
$wsExpo.Range("F37").Value = 50
$wsExpo.Range("F39").Value = 1208
$wsExpo.Range("F40").Value = 423
$wsExpo.Range("F43").Value = 230

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 7754
$wsShow.Range("F6").Value = 130
$wsShow.Range("F20").Value = 41
$wsShow.Range("F24").Value = 135
$wsShow.Range("F27").Value = 5
$wsShow.Range("F45").Value = 48
$wsShow.Range("F47").Value = 69
$wsShow.Range("F49").Value = 36

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 2407
$wsLocal.Range("F7").Value = 690
$wsLocal.Range("F9").Value = 9481
$wsLocal.Range("F11").Value = 194
$wsLocal.Range("F15").Value = 300
$wsLocal.Range("F16").Value = 2568
$wsLocal.Range("F17").Value = 272
$wsLocal.Range("F18").Value = 88
$wsLocal.Range("F19").Value = 559

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3603
$wsAll.Range("F4").Value = 2407
$wsAll.Range("F5").Value = 690
$wsAll.Range("F7").Value = 300
$wsAll.Range("F8").Value = 2568
$wsAll.Range("F9").Value = 272
$wsAll.Range("F10").Value = 680
$wsAll.Range("F11").Value = 7616
$wsAll.Range("F12").Value = 7828
$wsAll.Range("F13").Value = 4980
$wsAll.Range("F14").Value = 539
$wsAll.Range("G15").Value = 58.8
$wsAll.Range("F18").Value = 5355
$wsAll.Range("F21").Value = 88
$wsAll.Range("F22").Value = 1465
$wsAll.Range("F23").Value = 2071
$wsAll.Range("F24").Value = 130
$wsAll.Range("F25").Value = 559
$wsAll.Range("F30").Value = 198
$wsAll.Range("F31").Value = 1099
$wsAll.Range("F33").Value = 50
$wsAll.Range("F38").Value = 423
$wsAll.Range("F40").Value = 5
$wsAll.Range("F42").Value = 230
